$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update header row text (columns B..I get new labels) ---
$ws.Range("B1").Value = "ApplicationType"
$ws.Range("C1").Value = "FromTime"
$ws.Range("D1").Value = "ToTime"
$ws.Range("G1").Value = "Reason"
$ws.Range("H1").Value = "TotalDays"
$ws.Range("I1").Value = "TotalHours"

# --- Add two new header columns ---
$ws.Range("J1").Value = "StartDuration"
$ws.Range("K1").Value = "EndDuration"

# --- Add formatted data row (row 2) ---
$ws.Range("C2").NumberFormat = "mm:ss.0"
$ws.Range("D2").NumberFormat = "mm:ss.0"
$ws.Range("E2").NumberFormat = "mm-dd-yy"
$ws.Range("F2").NumberFormat = "mm-dd-yy"
$ws.Range("M2").NumberFormat = "mm:ss.0"

# --- Update column widths (best-fit-like) ---
$ws.Columns.Item(1).ColumnWidth = 6.7109375
$ws.Columns.Item(2).ColumnWidth = 15
$ws.Columns.Item(3).ColumnWidth = 9.7109375
$ws.Columns.Item(4).ColumnWidth = 7.28515625
$ws.Columns.Item(5).ColumnWidth = 9.7109375
$ws.Columns.Item(6).ColumnWidth = 7.28515625
$ws.Columns.Item(7).ColumnWidth = 7.5703125
$ws.Columns.Item(8).ColumnWidth = 9.5703125
$ws.Columns.Item(9).ColumnWidth = 10.5703125
$ws.Columns.Item(10).ColumnWidth = 12.7109375
$ws.Columns.Item(11).ColumnWidth = 11.85546875

# --- Update selection ---
$ws.Range("B2").Select()
